# Redo naming of mixres units
# - Remove the old "area_mixres" sheet (superseded data).
# - Rename "area_mixres_new" to "area_mixre", becoming the first sheet.
# - Remaining sheets (area_hires, area_lores, area_pop_sum) are left untouched.

$wb = $excel.ActiveWorkbook

$oldSheet = $wb.Worksheets.Item("area_mixres")
$oldSheet.Delete() | Out-Null

$newSheet = $wb.Worksheets.Item("area_mixres_new")
$newSheet.Name = "area_mixre"
